$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 12751
$ws.Range("I116").Value = 13668
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 13668
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = -10226
$ws.Range("N116").Value = -16884
$ws.Range("H131").Value = 6673.486
$ws.Range("I131").Value = 1667.3077
$ws.Range("K131").Value = 5001.9231
$ws.Range("M131").Value = 38.07690000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10788608
$ws.Range("I2").Value = 13215727
$ws.Range("J2").Value = 1410.1111
$ws.Range("K2").Value = 13215727
$ws.Range("L2").Value = 1410.1111
$ws.Range("M2").Value = -13215614
$ws.Range("N2").Value = -1636.1111
$ws.Range("H4").Value = 424.7143
$ws.Range("I4").Value = 144.53847
$ws.Range("J4").Value = 880
$ws.Range("K4").Value = 144.53847
$ws.Range("L4").Value = 880
$ws.Range("M4").Value = -28.53846999999999
$ws.Range("N4").Value = -1112
$ws.Range("H32").Value = 3574.9324
$ws.Range("I32").Value = 2681.6865
$ws.Range("K32").Value = 2681.6865
$ws.Range("M32").Value = -2394.6865
$ws.Range("H45").Value = 19441.428
$ws.Range("I45").Value = 22431.666
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 22431.666
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -22054.666
$ws.Range("N45").Value = -2254
$ws.Range("H61").Value = 15861.3
$ws.Range("I61").Value = 19769.166
$ws.Range("J61").Value = 9999.5
$ws.Range("K61").Value = 19769.166
$ws.Range("L61").Value = 9999.5
$ws.Range("M61").Value = -19557.166
$ws.Range("N61").Value = -10423.5
$ws.Range("H74").Value = 7191.2085
$ws.Range("I74").Value = 7191.2085
$ws.Range("K74").Value = 7191.2085
$ws.Range("M74").Value = -6317.2085
$ws.Range("H77").Value = 7191.2085
$ws.Range("I77").Value = 7191.2085
$ws.Range("K77").Value = 35956.0425
$ws.Range("M77").Value = -31588.0425
$ws.Range("H102").Value = 5666.722
$ws.Range("I102").Value = 5418.4116
$ws.Range("J102").Value = 9888
$ws.Range("K102").Value = 5418.4116
$ws.Range("L102").Value = 9888
$ws.Range("M102").Value = -3796.4116
$ws.Range("N102").Value = -13132
$ws.Range("H110").Value = 1571.9048
$ws.Range("J110").Value = 2308
$ws.Range("L110").Value = 2308
$ws.Range("N110").Value = -6398
$ws.Range("H116").Value = 10788608
$ws.Range("I116").Value = 13215727
$ws.Range("J116").Value = 1410.1111
$ws.Range("K116").Value = 13215727
$ws.Range("L116").Value = 1410.1111
$ws.Range("M116").Value = -13213433
$ws.Range("N116").Value = -5998.1111
$ws.Range("H132").Value = 3773.75
$ws.Range("I132").Value = 3642.7222
$ws.Range("J132").Value = 4166.8335
$ws.Range("K132").Value = 10928.1666
$ws.Range("L132").Value = 12500.5005
$ws.Range("M132").Value = -8398.1666
$ws.Range("N132").Value = -17560.5005
$ws.Range("H136").Value = 15861.3
$ws.Range("I136").Value = 19769.166
$ws.Range("J136").Value = 9999.5
$ws.Range("K136").Value = 59307.49800000001
$ws.Range("L136").Value = 29998.5
$ws.Range("M136").Value = -56757.49800000001
$ws.Range("N136").Value = -35098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10788608
$ws.Range("I3").Value = 13215727
$ws.Range("J3").Value = 1410.1111
$ws.Range("K3").Value = 13215727
$ws.Range("L3").Value = 1410.1111
$ws.Range("M3").Value = -13215613
$ws.Range("N3").Value = -1638.1111
$ws.Range("H20").Value = 2502.6924
$ws.Range("I20").Value = 1867.9474
$ws.Range("J20").Value = 4225.5713
$ws.Range("K20").Value = 1867.9474
$ws.Range("L20").Value = 4225.5713
$ws.Range("M20").Value = -1620.9474
$ws.Range("N20").Value = -4719.5713
$ws.Range("H86").Value = 2352.4285
$ws.Range("I86").Value = 2318.1428
$ws.Range("K86").Value = 2318.1428
$ws.Range("M86").Value = -1195.1428
$ws.Range("H89").Value = 2352.4285
$ws.Range("I89").Value = 2318.1428
$ws.Range("K89").Value = 11590.714
$ws.Range("M89").Value = -5974.714
$ws.Range("H102").Value = 42178.285
$ws.Range("I102").Value = 64099.2
$ws.Range("K102").Value = 64099.2
$ws.Range("M102").Value = -60854.2
$ws.Range("H105").Value = 4004.7778
$ws.Range("I105").Value = 2395.9
$ws.Range("K105").Value = 2395.9
$ws.Range("M105").Value = -648.9000000000001
$ws.Range("H107").Value = 4891.9
$ws.Range("I107").Value = 4614.6
$ws.Range("J107").Value = 5723.8
$ws.Range("K107").Value = 4614.6
$ws.Range("L107").Value = 5723.8
$ws.Range("M107").Value = -2694.6
$ws.Range("N107").Value = -9563.799999999999
$ws.Range("H134").Value = 4833.2954
$ws.Range("I134").Value = 4767.5854
$ws.Range("K134").Value = 14302.7562
$ws.Range("M134").Value = -11767.7562

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1691.9524
$ws.Range("I16").Value = 1461.7693
$ws.Range("J16").Value = 2066
$ws.Range("K16").Value = 1461.7693
$ws.Range("L16").Value = 2066
$ws.Range("M16").Value = -1174.7693
$ws.Range("N16").Value = -2640
$ws.Range("H31").Value = 2495.8655
$ws.Range("I31").Value = 2240.9048
$ws.Range("J31").Value = 3566.7
$ws.Range("K31").Value = 2240.9048
$ws.Range("L31").Value = 3566.7
$ws.Range("M31").Value = -1945.9048
$ws.Range("N31").Value = -4156.7
$ws.Range("H34").Value = 2495.8655
$ws.Range("I34").Value = 2240.9048
$ws.Range("J34").Value = 3566.7
$ws.Range("K34").Value = 2240.9048
$ws.Range("L34").Value = 3566.7
$ws.Range("M34").Value = -2038.9048
$ws.Range("N34").Value = -3970.7
$ws.Range("H58").Value = 7491.185
$ws.Range("I58").Value = 9333.817999999999
$ws.Range("K58").Value = 9333.817999999999
$ws.Range("M58").Value = -9130.817999999999
$ws.Range("H105").Value = 1647
$ws.Range("I105").Value = 1662.7368
$ws.Range("K105").Value = 1662.7368
$ws.Range("M105").Value = 84.2632000000001
$ws.Range("H107").Value = 1115.6
$ws.Range("I107").Value = 716.2222
$ws.Range("K107").Value = 716.2222
$ws.Range("M107").Value = 1203.7778
$ws.Range("H113").Value = 1691.9524
$ws.Range("I113").Value = 1461.7693
$ws.Range("J113").Value = 2066
$ws.Range("K113").Value = 1461.7693
$ws.Range("L113").Value = 2066
$ws.Range("M113").Value = 708.2307000000001
$ws.Range("N113").Value = -6406
$ws.Range("H134").Value = 31562.2
$ws.Range("I134").Value = 36452.75
$ws.Range("J134").Value = 12000
$ws.Range("K134").Value = 109358.25
$ws.Range("L134").Value = 36000
$ws.Range("M134").Value = -106823.25
$ws.Range("N134").Value = -41070
$ws.Range("H136").Value = 7491.185
$ws.Range("I136").Value = 9333.817999999999
$ws.Range("K136").Value = 28001.454
$ws.Range("M136").Value = -25451.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 366.7143
$ws.Range("I107").Value = 261.8
$ws.Range("K107").Value = 785.4000000000001
$ws.Range("M107").Value = 1134.6
$ws.Range("H129").Value = 14941520
$ws.Range("I129").Value = 17857858
$ws.Range("J129").Value = 12349220
$ws.Range("K129").Value = 53573574
$ws.Range("L129").Value = 37047660
$ws.Range("M129").Value = -53568574
$ws.Range("N129").Value = -37057660
$ws.Range("H131").Value = 5263262.5
$ws.Range("I131").Value = 7520872
$ws.Range("J131").Value = 3973199.8
$ws.Range("K131").Value = 22562616
$ws.Range("L131").Value = 11919599.4
$ws.Range("M131").Value = -22557576
$ws.Range("N131").Value = -11929679.4
$ws.Range("H139").Value = 2574.3667
$ws.Range("I139").Value = 1176.5
$ws.Range("J139").Value = 8165.8335
$ws.Range("K139").Value = 3529.5
$ws.Range("L139").Value = 24497.5005
$ws.Range("M139").Value = 1610.5
$ws.Range("N139").Value = -34777.50049999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 23002564
$ws.Range("I80").Value = 30668352
$ws.Range("J80").Value = 5199.2
$ws.Range("K80").Value = 30668352
$ws.Range("L80").Value = 5199.2
$ws.Range("M80").Value = -30667354
$ws.Range("N80").Value = -7195.2
$ws.Range("H83").Value = 23002564
$ws.Range("I83").Value = 30668352
$ws.Range("J83").Value = 5199.2
$ws.Range("K83").Value = 153341760
$ws.Range("L83").Value = 25996
$ws.Range("M83").Value = -153336768
$ws.Range("N83").Value = -35980
$ws.Range("H102").Value = 5480.773
$ws.Range("I102").Value = 8289.6
$ws.Range("K102").Value = 8289.6
$ws.Range("M102").Value = -6667.6
$ws.Range("H113").Value = 3470.25
$ws.Range("I113").Value = 2440.5
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 2440.5
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -270.5
$ws.Range("N113").Value = -8840
$ws.Range("H122").Value = 2826.7778
$ws.Range("I122").Value = 3054.6316
$ws.Range("J122").Value = 2285.625
$ws.Range("K122").Value = 9163.8948
$ws.Range("L122").Value = 6856.875
$ws.Range("M122").Value = -6713.8948
$ws.Range("N122").Value = -11756.875
$ws.Range("H126").Value = 7486.8
$ws.Range("I126").Value = 5812.3335
$ws.Range("J126").Value = 9998.5
$ws.Range("K126").Value = 17437.0005
$ws.Range("L126").Value = 29995.5
$ws.Range("M126").Value = -14967.0005
$ws.Range("N126").Value = -34935.5
$ws.Range("H132").Value = 7567.75
$ws.Range("I132").Value = 8824.556
$ws.Range("J132").Value = 6539.4546
$ws.Range("K132").Value = 26473.668
$ws.Range("L132").Value = 19618.3638
$ws.Range("M132").Value = -23943.668
$ws.Range("N132").Value = -24678.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 947.73334
$ws.Range("I93").Value = 561
$ws.Range("K93").Value = 561
$ws.Range("M93").Value = 687

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 129999
$ws.Range("J82").Value = 129999
$ws.Range("L82").Value = 129999
$ws.Range("N82").Value = -130765
$ws.Range("H85").Value = 129999
$ws.Range("J85").Value = 129999
$ws.Range("L85").Value = 129999
$ws.Range("N85").Value = -132651
$ws.Range("H126").Value = 7173.905
$ws.Range("J126").Value = 10665.333
$ws.Range("L126").Value = 31995.999
$ws.Range("N126").Value = -36935.999
$ws.Range("H133").Value = 37100
$ws.Range("J133").Value = 37100
$ws.Range("L133").Value = 37100
$ws.Range("N133").Value = -47220
